$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.187.65"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "2.427.80"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'490.37"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("E6").Value = "  +3.84%  "
$ws.Range("D7").Value = "'0.619"
$ws.Range("E7").Value = "  +20.25%  "
$ws.Range("D8").Value = "'0.997"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Value = "2.447.76"
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("D11").Value = "'5.67"
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D14").Value = "2.858.55"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "57.328.11"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").Value = "'20.86"
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("D17").Value = "'0.0000134"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").Value = "2.443.29"
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("D19").Value = "'4.77"
$ws.Range("E19").Value = "  +5.60%  "
$ws.Range("D20").Value = "'328.49"
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("D21").Value = "'9.99"
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'5.93"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Value = "'58.37"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").Value = "2.535.94"
$ws.Range("E28").Value = "  -2.41%  "
$ws.Range("D29").Value = "'7.32"
$ws.Range("E29").Value = "  -4.06%  "
$ws.Range("D30").Value = "0.0₃0794"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("D33").Value = "'149.44"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "'1.52"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").Value = "'5.34"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("D38").Value = "'0.861"
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("E39").Value = "  +11.34%  "
$ws.Range("D40").Value = "'34.25"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("D41").Value = "'1.38"
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("D42").Value = "'3.53"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "'0.598"
$ws.Range("E44").Value = "  -2.05%  "
$ws.Range("D45").Value = "'0.0537"
$ws.Range("E45").Value = "  -3.69%  "
$ws.Range("D46").Value = "'267.28"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "'4.69"
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.862.64"
$ws.Range("E51").Value = "  -1.56%  "
